$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Row 1 (header row) ---
# New header labels introduced first in the exact order needed so that the
# shared-string table grows in the same order as the target workbook.
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"

# --- Row 2 / Row 3 values (this is what drives most of the new shared strings) ---
$ws.Range("B2").Value = "土地抵押權"
$ws.Range("B3").Value = "房貸"

$ws.Range("D2").Value = "華泰銀行臺北市長安東路"
$ws.Range("D3").Value = "合作金庫臺北市和平東路"

$ws.Range("F2").Value = "99年11月05日"
$ws.Range("F3").Value = "99年06月30日"

$ws.Range("G2").Value = "抵押借款"
$ws.Range("G3").Value = "房貸"

$ws.Range("H2").Value = "debt"
$ws.Range("H3").Value = "debt"

# --- Finish header row with the standard trailing columns ---
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Remaining data cells for rows 2 and 3 ---
$ws.Range("C2").Value = "薛凌"
$ws.Range("C3").Value = "陳勝宏"

$ws.Range("I2").Value = "normal"
$ws.Range("I3").Value = "normal"

$ws.Range("J2").Value = "2011-11-14"
$ws.Range("J3").Value = "2011-11-14"

$ws.Range("K2").Value = "薛凌"
$ws.Range("K3").Value = "薛凌"

$ws.Range("L2").Value = 1384
$ws.Range("L3").Value = 1384

$ws.Range("M2").Value = "tmp7f9c1"
$ws.Range("M3").Value = "tmp7f9c1"

$ws.Range("N2").Value = 219
$ws.Range("N3").Value = 220

# Header row style (bold, bordered) to match the rest of the sheet
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
